# Update "想去人数" (F column) counts across sheets, as published at 456a3b4.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 526
$ws1.Range("F6").Value  = 493
$ws1.Range("F7").Value  = 95
$ws1.Range("F8").Value  = 108
$ws1.Range("F10").Value = 6566
$ws1.Range("F12").Value = 358
$ws1.Range("F13").Value = 2781
$ws1.Range("F15").Value = 292
$ws1.Range("F17").Value = 517

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 12

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 12
$ws4.Range("F5").Value  = 526
$ws4.Range("F8").Value  = 493
$ws4.Range("F9").Value  = 95
$ws4.Range("F10").Value = 108
$ws4.Range("F13").Value = 6566
$ws4.Range("F16").Value = 358
$ws4.Range("F17").Value = 2781
$ws4.Range("F19").Value = 292
$ws4.Range("F21").Value = 517
